# Release v0.1.0-beta: Fix validation errors and update canonical URL
#
# Updates the "Metadata" sheet (Version, Status, Date, Description) and the
# "Elements" sheet (Definition of the root Extension element) to match the
# new release.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 0.1.0
$meta.Range("B3").Value = "0.1.0"

# Status: active -> draft
$meta.Range("B6").Value = "draft"

# Date: refreshed publish timestamp
$meta.Range("B8").Value = "2025-12-26T14:13:58+00:00"

# Description: was empty, now populated
$meta.Range("B11").Value = "Extension to link goal evaluation observations to the patient goals being evaluated. Enables tracking of goal progress and outcomes over time."

# ---- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The "Definition" column (M) for the root Extension element (row 2) picked
# up the StructureDefinition's new description text instead of the generic
# placeholder "An Extension".
$elements.Range("M2").Value = "Extension to link goal evaluation observations to the patient goals being evaluated. Enables tracking of goal progress and outcomes over time."
